# Mark attendance results for each session row: set the appropriate
# "Real"/"Total Attendance Count" (D/E) cells to 1 for sessions the
# student attended, and the "Absent" (H) cell to 1 for sessions missed.
# Row 3 also gets its "Invalid" (G) flag set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToMark = @(
    "G3", "H3",
    "D4", "E4",
    "D5", "E5",
    "D6", "E6",
    "H7",
    "H8",
    "D9", "E9",
    "H10",
    "D11", "E11",
    "D12", "E12",
    "D13", "E13",
    "D14", "E14",
    "H15",
    "H16",
    "D17", "E17",
    "H18"
)

foreach ($addr in $cellsToMark) {
    $ws.Range($addr).Value = 1
}
